$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header text for column O (shared string used by O1 changes meaning)
$ws.Range("O1").Value = "寻址方式（可选）"

# O2 was boolean TRUE -> now text "静态"
$ws.Range("O2").Value = "静态"

# O3 new cell with text "动态"
$ws.Range("O3").Value = "动态"

# Update selection to O3 to match the saved selection state
$ws.Range("O3").Select()
